$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Vision" subsystem entry at row 10 (PORT 9)
$ws.Range("B10").Value = "Vision"
$ws.Range("C10").Value = "vision"
$ws.Range("D10").Value = "Vision Sensor"

# Update selection to match the new active range
$ws.Range("B10:D10").Select()
